$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '88.922.02'
$ws.Range('E2').Value = '  +9.43%  '

# Row 3
$ws.Range('D3').Value = '3.369.87'
$ws.Range('E3').Value = '  +7.09%  '

# Row 4
$ws.Range('E4').Value = '  +0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '221.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.72%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '650.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.00%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.408'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +43.73%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.619'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.87%  '

# Row 10
$ws.Range('D10').Value = '3.365.79'
$ws.Range('E10').Value = '  +7.04%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.642'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +10.98%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000289'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +14.81%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.98'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +17.78%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.169'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.42%  '

# Row 15
$ws.Range('D15').Value = '4.003.42'
$ws.Range('E15').Value = '  +7.52%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.56'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.21%  '

# Row 17
$ws.Range('D17').Value = '88.381.37'
$ws.Range('E17').Value = '  +9.07%  '

# Row 18
$ws.Range('D18').Value = '3.365.36'
$ws.Range('E18').Value = '  +7.21%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.73%  '

# Row 20
$ws.Range('E20').Value = '  -0.95%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.98%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '458.22'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.34%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +10.59%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.15%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.60'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.81%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +20.05%  '

# Row 27
$ws.Range('D27').Value = '3.549.66'
$ws.Range('E27').Value = '  +7.43%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000143'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +18.19%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '79.64'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.16%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.200'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +46.54%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.19%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.51'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.28%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '594.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.38%  '

# Row 34
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.997'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.25%  '

# Row 35
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.43%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.10%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.34'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +20.54%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.145'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.34%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.61'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.25%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.432'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.97%  '

# Row 41
$ws.Range('E41').Value = '  +6.44%  '

# Row 42
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.82'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.30%  '

# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.11%  '

# Row 44
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.16'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.17%  '

# Row 45
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '158.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.06%  '

# Row 46
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.46'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +10.57%  '

# Row 47
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.02%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '189.09'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.17%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.83%  '

# Row 50
$ws.Range('E50').Value = '  +7.92%  '

# Row 51
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.790'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.12%  '
